# Generate Report for Handback
#
# - Overview/zh-cn/de-de "Status" cells flip from "Ready for handoff" to
#   "Handed back: in sync with en-US" (all cells share the one string).
# - zh-cn & de-de sheets each gain "Latest Target File" (E) / "Latest
#   Handback File" (F) hyperlink cells for rows 2 & 3, and the "Latest
#   Handback DateTime" (G) for rows 2 & 3 is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# ---- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

$mdUrl_26bed90c = "https://github.com/OpenLocalizationTest/oltest/blob/84068b1bb824c7cd2eda295be60f6afc3c233c16/e2e/26bed90c-d4f9-4756-9a78-655843e4d9e3.md"
$mdName_26bed90c = "26bed90c-d4f9-4756-9a78-655843e4d9e3.md"

$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0154f79630b25ae202d4bd8346ec7425816e531c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/26bed90c-d4f9-4756-9a78-655843e4d9e3.7f5c799414885d0543798ee12ab51b45f4691c63.zh-cn.xlf"
$zhXlfName = "26bed90c-d4f9-4756-9a78-655843e4d9e3.7f5c799414885d0543798ee12ab51b45f4691c63.zh-cn.xlf"

$zh.Hyperlinks.Add($zh.Range("E2"), $mdUrl_26bed90c, "", "", $mdName_26bed90c)
$zh.Hyperlinks.Add($zh.Range("F2"), $zhXlfUrl, "", "", $zhXlfName)
$zh.Hyperlinks.Add($zh.Range("E3"), $mdUrl_26bed90c, "", "", $mdName_26bed90c)
$zh.Hyperlinks.Add($zh.Range("F3"), $zhXlfUrl, "", "", $zhXlfName)

$zh.Range("E2").Font.Underline = 2
$zh.Range("E2").Font.Color = 15570276
$zh.Range("F2").Font.Underline = 2
$zh.Range("F2").Font.Color = 15570276
$zh.Range("E3").Font.Underline = 2
$zh.Range("E3").Font.Color = 15570276
$zh.Range("F3").Font.Underline = 2
$zh.Range("F3").Font.Color = 15570276

$zh.Range("G2").Value = "2016-01-25 08:47:42"
$zh.Range("G3").Value = "2016-01-25 08:47:42"

# ---- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ab7cf0cc310db707fe54b2eab6b53abd7e9bb4e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/26bed90c-d4f9-4756-9a78-655843e4d9e3.7f5c799414885d0543798ee12ab51b45f4691c63.de-de.xlf"
$deXlfName = "26bed90c-d4f9-4756-9a78-655843e4d9e3.7f5c799414885d0543798ee12ab51b45f4691c63.de-de.xlf"

$de.Hyperlinks.Add($de.Range("E2"), $mdUrl_26bed90c, "", "", $mdName_26bed90c)
$de.Hyperlinks.Add($de.Range("F2"), $deXlfUrl, "", "", $deXlfName)
$de.Hyperlinks.Add($de.Range("E3"), $mdUrl_26bed90c, "", "", $mdName_26bed90c)
$de.Hyperlinks.Add($de.Range("F3"), $deXlfUrl, "", "", $deXlfName)

$de.Range("E2").Font.Underline = 2
$de.Range("E2").Font.Color = 15570276
$de.Range("F2").Font.Underline = 2
$de.Range("F2").Font.Color = 15570276
$de.Range("E3").Font.Underline = 2
$de.Range("E3").Font.Color = 15570276
$de.Range("F3").Font.Underline = 2
$de.Range("F3").Font.Color = 15570276

$de.Range("G2").Value = "2016-01-25 08:48:00"
$de.Range("G3").Value = "2016-01-25 08:48:00"
